# Update cryptos list values (Price column D, Volume(1h) column E)
# per refreshed data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "27.184.94"
$ws.Range('E2').Value = "  -0.91%  "
$ws.Range('D3').Value = "1.646.99"
$ws.Range('E3').Value = "  -1.27%  "
$ws.Range('E4').Value = "  +0.16%  "
$ws.Range('D5').Value = "'218.55"
$ws.Range('D6').Value = "'0.510"
$ws.Range('E6').Value = "  +0.98%  "
$ws.Range('E7').Value = "  +0.16%  "
$ws.Range('E8').Value = "  +0.37%  "
$ws.Range('D9').Value = "'0.0628"
$ws.Range('E9').Value = "  -0.23%  "
$ws.Range('D10').Value = "'20.04"
$ws.Range('E10').Value = "  +0.39%  "
$ws.Range('E11').Value = "  -0.40%  "
$ws.Range('D12').Value = "1.877.79"
$ws.Range('E12').Value = "  -1.16%  "
$ws.Range('D13').Value = "1.645.07"
$ws.Range('E13').Value = "  -1.41%  "
$ws.Range('E14').Value = "  -2.19%  "
$ws.Range('E15').Value = "  +0.45%  "
$ws.Range('D16').Value = "'67.40"
$ws.Range('E16').Value = "  +0.02%  "
$ws.Range('D17').Value = "27.187.48"
$ws.Range('E17').Value = "  -0.74%  "
$ws.Range('E18').Value = "  +0.05%  "
$ws.Range('D19').Value = "'219.63"
$ws.Range('E19').Value = "  -2.43%  "
$ws.Range('E20').Value = "  +0.08%  "
$ws.Range('D21').Value = "'6.78"
$ws.Range('E21').Value = "  -0.35%  "
$ws.Range('E22').Value = "  -0.41%  "
$ws.Range('D23').Value = "'2.50"
$ws.Range('E23').Value = "  +1.77%  "
$ws.Range('E24').Value = "  -1.12%  "
$ws.Range('D25').Value = "'148.33"
$ws.Range('E25').Value = "  +0.47%  "
$ws.Range('E26').Value = "  +0.17%  "
$ws.Range('D27').Value = "'7.41"
$ws.Range('E27').Value = "  -1.29%  "
$ws.Range('E28').Value = "  -0.47%  "
$ws.Range('D29').Value = "'15.79"
$ws.Range('E29').Value = "  -1.65%  "
$ws.Range('E30').Value = "  -1.72%  "
$ws.Range('E31').Value = "  -1.06%  "
$ws.Range('D32').Value = "'3.36"
$ws.Range('E32').Value = "  -0.98%  "
$ws.Range('E33').Value = "  +0.73%  "
$ws.Range('E34').Value = "  +0.78%  "
$ws.Range('D35').Value = "1.269.53"
$ws.Range('E35').Value = "  -0.58%  "
$ws.Range('D36').Value = "'2.46"
$ws.Range('E36').Value = "  +0.63%  "
$ws.Range('E37').Value = "  +0.28%  "
$ws.Range('E38').Value = "  +0.59%  "
$ws.Range('E39').Value = "  +1.16%  "
$ws.Range('E40').Value = "  +0.11%  "
$ws.Range('D41').Value = "'0.810"
$ws.Range('E41').Value = "  -0.41%  "
$ws.Range('E42').Value = "  +4.38%  "
$ws.Range('E43').Value = "  -0.32%  "
$ws.Range('D44').Value = "1.788.68"
$ws.Range('E44').Value = "  -1.21%  "
$ws.Range('D45').Value = "'62.51"
$ws.Range('E45').Value = "  +0.24%  "
$ws.Range('D46').Value = "'92.37"
$ws.Range('E46').Value = "  -0.39%  "
$ws.Range('D47').Value = "'1.60"
$ws.Range('E47').Value = "  -2.18%  "
$ws.Range('D48').Value = "'0.0513"
$ws.Range('E48').Value = "  -0.95%  "
$ws.Range('D49').Value = "'7.69"
$ws.Range('E49').Value = "  -0.33%  "
$ws.Range('E50').Value = "  -1.34%  "
$ws.Range('D51').Value = "'0.405"
$ws.Range('E51').Value = "  -0.75%  "
